# error solve ifrs list
# Update the financial figures for the 5 existing data rows (rows 2-6) and
# remove the (erroneous/duplicated) data that had been entered for rows 7-9,
# leaving only the identifying columns (A, B, C) for those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 ----
$ws.Range("D2").Value  = 4284
$ws.Range("E2").Value  = 316
$ws.Range("F2").Value  = 316
$ws.Range("G2").Value  = 287
$ws.Range("H2").Value  = 214
$ws.Range("I2").Value  = 195
$ws.Range("J2").Value  = 19
$ws.Range("K2").Value  = 4000
$ws.Range("L2").Value  = 1568
$ws.Range("M2").Value  = 2432
$ws.Range("N2").Value  = 1899
$ws.Range("O2").Value  = 533
$ws.Range("P2").Value  = 343
$ws.Range("Q2").Value  = 251
$ws.Range("R2").Value  = -42
$ws.Range("S2").Value  = -124
$ws.Range("T2").Value  = 79
$ws.Range("U2").Value  = 172
$ws.Range("V2").Value  = 578
$ws.Range("W2").Value  = 7.37
$ws.Range("X2").Value  = 4.99
$ws.Range("Y2").Value  = 10.68
$ws.Range("Z2").Value  = 5.35
$ws.Range("AA2").Value = 64.45999999999999
$ws.Range("AB2").Value = 506.02
$ws.Range("AC2").Value = 284
$ws.Range("AD2").Value = 19.51
$ws.Range("AE2").Value = 3208
$ws.Range("AF2").Value = 1.73
$ws.Range("AG2").Value = 68
$ws.Range("AH2").Value = 1.23
$ws.Range("AI2").Value = 20.64
$ws.Range("AJ2").Value = 68560000

# ---- Row 3 ----
$ws.Range("D3").Value  = 4428
$ws.Range("E3").Value  = 399
$ws.Range("F3").Value  = 399
$ws.Range("G3").Value  = 399
$ws.Range("H3").Value  = 303
$ws.Range("I3").Value  = 279
$ws.Range("J3").Value  = 24
$ws.Range("K3").Value  = 4295
$ws.Range("L3").Value  = 1621
$ws.Range("M3").Value  = 2674
$ws.Range("N3").Value  = 2121
$ws.Range("O3").Value  = 553
$ws.Range("P3").Value  = 343
$ws.Range("Q3").Value  = 408
$ws.Range("R3").Value  = -115
$ws.Range("S3").Value  = -54
$ws.Range("T3").Value  = 49
$ws.Range("U3").Value  = 359
$ws.Range("V3").Value  = 567
$ws.Range("W3").Value  = 9.01
$ws.Range("X3").Value  = 6.84
$ws.Range("Y3").Value  = 13.9
$ws.Range("Z3").Value  = 7.3
$ws.Range("AA3").Value = 60.64
$ws.Range("AB3").Value = 570.66
$ws.Range("AC3").Value = 407
$ws.Range("AD3").Value = 20.87
$ws.Range("AE3").Value = 3583
$ws.Range("AF3").Value = 2.37
$ws.Range("AG3").Value = 65
$ws.Range("AH3").Value = 0.76
$ws.Range("AI3").Value = 13.78
$ws.Range("AJ3").Value = 68560000

# ---- Row 4 ----
$ws.Range("D4").Value  = 4438
$ws.Range("E4").Value  = 352
$ws.Range("F4").Value  = 352
$ws.Range("G4").Value  = 321
$ws.Range("H4").Value  = 244
$ws.Range("I4").Value  = 225
$ws.Range("J4").Value  = 19
$ws.Range("K4").Value  = 4584
$ws.Range("L4").Value  = 1781
$ws.Range("M4").Value  = 2803
$ws.Range("N4").Value  = 2235
$ws.Range("O4").Value  = 567
$ws.Range("P4").Value  = 343
$ws.Range("Q4").Value  = 386
$ws.Range("R4").Value  = -220
$ws.Range("S4").Value  = -94
$ws.Range("T4").Value  = 355
$ws.Range("U4").Value  = 31
$ws.Range("V4").Value  = 562
$ws.Range("W4").Value  = 7.92
$ws.Range("X4").Value  = 5.5
$ws.Range("Y4").Value  = 10.31
$ws.Range("Z4").Value  = 5.5
$ws.Range("AA4").Value = 63.56
$ws.Range("AB4").Value = 618.72
$ws.Range("AC4").Value = 328
$ws.Range("AD4").Value = 14.23
$ws.Range("AE4").Value = 3840
$ws.Range("AF4").Value = 1.21
$ws.Range("AG4").Value = 35
$ws.Range("AH4").Value = 0.75
$ws.Range("AI4").Value = 9.07
$ws.Range("AJ4").Value = 68560000

# ---- Row 5 ----
$ws.Range("D5").Value  = 4147
$ws.Range("E5").Value  = 216
$ws.Range("F5").Value  = 216
$ws.Range("G5").Value  = 185
$ws.Range("H5").Value  = 134
$ws.Range("I5").Value  = 108
$ws.Range("J5").Value  = 25
$ws.Range("K5").Value  = 4516
$ws.Range("L5").Value  = 1598
$ws.Range("M5").Value  = 2918
$ws.Range("N5").Value  = 2327
$ws.Range("O5").Value  = 592
$ws.Range("P5").Value  = 343
$ws.Range("Q5").Value  = 169
$ws.Range("R5").Value  = -166
$ws.Range("S5").Value  = -96
$ws.Range("T5").Value  = 222
$ws.Range("U5").Value  = -53
$ws.Range("V5").Value  = 487
$ws.Range("W5").Value  = 5.22
$ws.Range("X5").Value  = 3.22
$ws.Range("Y5").Value  = 4.74
$ws.Range("Z5").Value  = 2.94
$ws.Range("AA5").Value = 54.76
$ws.Range("AB5").Value = 644.66
$ws.Range("AC5").Value = 158
$ws.Range("AD5").Value = 25.37
$ws.Range("AE5").Value = 3997
$ws.Range("AF5").Value = 1
$ws.Range("AG5").Value = 25
$ws.Range("AH5").Value = 0.63
$ws.Range("AI5").Value = 13.46
$ws.Range("AJ5").Value = 68560000

# ---- Row 6 (note: J6, AG6 and AH6 are not populated in the new data) ----
$ws.Range("D6").Value  = 4324
$ws.Range("E6").Value  = 99
$ws.Range("F6").Value  = 99
$ws.Range("G6").Value  = 119
$ws.Range("H6").Value  = 84
$ws.Range("I6").Value  = 79
$ws.Range("K6").Value  = 4746
$ws.Range("L6").Value  = 1789
$ws.Range("M6").Value  = 2957
$ws.Range("N6").Value  = 2374
$ws.Range("P6").Value  = 343
$ws.Range("Q6").Value  = -49
$ws.Range("R6").Value  = -439
$ws.Range("S6").Value  = 281
$ws.Range("T6").Value  = 401
$ws.Range("U6").Value  = -450
$ws.Range("V6").Value  = 808
$ws.Range("W6").Value  = 2.29
$ws.Range("X6").Value  = 1.95
$ws.Range("Y6").Value  = 3.37
$ws.Range("Z6").Value  = 1.82
$ws.Range("AA6").Value = 60.48
$ws.Range("AB6").Value = 660.65
$ws.Range("AC6").Value = 115
$ws.Range("AD6").Value = 21.23
$ws.Range("AE6").Value = 4078
$ws.Range("AF6").Value = 0.6
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()
$ws.Range("AI6").Value = 11.04
$ws.Range("AJ6").Value = 68560000

# ---- Rows 7-9: the financial data columns (D:AI) were mistakenly filled in
#      and must be cleared, keeping only the A/B/C identifying columns ----
$ws.Range("D7:AI7").ClearContents()
$ws.Range("D8:AI8").ClearContents()
$ws.Range("D9:AI9").ClearContents()
